$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pour la prochaine fois")
$ws.Activate()

# Insert a new row at row 3, pushing the existing rows (3..21) down to (4..22).
$ws.Rows.Item(3).Insert()

# The task "Guncontroller voir la balle partir au tir" (still sitting in B2) is
# done -> clear its text and mark the row with a yellow highlight.
$ws.Range("B2").ClearContents()
$ws.Range("B2").Interior.Color = 65535

# New task dropped into the blank row that opened up between the old B6
# ("POV d'un joueur bug desfois...", now B7) and B9 ("joueur enemie...", now
# B10). Entered first so it claims the earlier shared-string slot.
$ws.Range("B9").Value = "Comprendre pq update function don't works for enemy entites"

# New task added right under the cleared/highlighted row.
$ws.Range("B3").Value = "Clique gauche viser"

# Match the author's final selection.
$ws.Range("B3").Select()
